$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# BOM part swaps (Buck capacitor changes):
#   Row 5: UUD1H150MCL1GS (15 uF Cap)   -> TR3E156M050C0300
#   Row 6: UWP1HR47MCL1GB (0.47uF Cap)  -> 293D474X9050B2TE3
#   Row 9: UCZ1J181MNJ1MS (180uF Cap)   -> TPSC107K010R0200 (100uF Cap)
# ---------------------------------------------------------------------------

# Row 5 - TR3E156M050C0300 (100uF->15uF cap slot), qty 6->5, price 0.48->2.32
$ws.Range("A5").Value = "TR3E156M050C0300"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 2.32

# Row 6 - 293D474X9050B2TE3, price 0.68->0.69
$ws.Range("A6").Value = "293D474X9050B2TE3"
$ws.Range("D6").Value = 0.69

# Row 9 - TPSC107K010R0200 / 100uF Cap, qty 4->6, price 1.87->0.82
$ws.Range("A9").Value = "TPSC107K010R0200"
$ws.Range("B9").Value = "100uF Cap"
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 0.82

# ---------------------------------------------------------------------------
# Rebuild the Digikey hyperlinks in column G (the underlying COM layer only
# supports clearing *all* hyperlinks on the sheet at once, so we capture the
# existing targets first, patch the three that changed, then recreate all of
# them in their original order).
# ---------------------------------------------------------------------------
$links = [ordered]@{
    "G2"  = "https://www.digikey.com/en/products/detail/stmicroelectronics/STM32F207VGT6TR/4357621"
    "G3"  = "https://www.digikey.com/en/products/detail/texas-instruments/LM2678S-5.0%2FNOPB/363825?utm_adgroup=General&utm_source=google&utm_medium=cpc&utm_campaign=PMax%20Shopping_Product_Zombie%20SKUs&utm_term=&utm_content=General&utm_id=go_cmp-17815035045_adg-_ad-__dev-c_ext-_prd-363825_sig-CjwKCAiA5L2tBhBTEiwAdSxJX2jmx8jM-JlHhV04F58rlCzi0KZgwJl8jmcjRGNCM7uSaMTsq63izRoCBfYQAvD_BwE&gad_source=1&gclid=CjwKCAiA5L2tBhBTEiwAdSxJX2jmx8jM-JlHhV04F58rlCzi0KZgwJl8jmcjRGNCM7uSaMTsq63izRoCBfYQAvD_BwE"
    "G4"  = "https://www.digikey.com/en/products/detail/texas-instruments/LM2678S-3-3-NOPB/366918"
    "G5"  = "https://www.digikey.com/en/products/detail/vishay-sprague/TR3E156M050C0300/2259941"
    "G6"  = "https://www.digikey.com/en/products/detail/vishay-sprague/293D474X9050B2TE3/1578900"
    "G7"  = "https://www.digikey.com/en/products/detail/vishay-general-semiconductor-diodes-division/VS-6TQ045S-M3/5426222"
    "G8"  = "https://www.digikey.com/en/products/detail/bourns-inc/2205-H-RC/775358"
    "G9"  = "https://www.digikey.com/en/products/detail/kyocera-avx/TPSC107K010R0200/946529"
    "G10" = "https://www.digikey.com/en/products/detail/murata-electronics/GRM2195C1H103JA01D/586788"
    "G11" = "https://www.digikey.com/en/products/detail/texas-instruments/SN75ALS174ADWR/1593485"
    "G12" = "https://www.digikey.com/en/products/detail/texas-instruments/SN65LBC175AD/380303"
    "G13" = "https://www.digikey.com/en/products/detail/cui-devices/PJ-202AH/408450"
    "G16" = "https://www.digikey.com/en/products/detail/texas-instruments/TPS63700DRCR/1672393"
    "G17" = "https://www.digikey.com/en/products/detail/vishay-general-semiconductor-diodes-division/SL03-GS18/4871689"
    "G18" = "https://www.digikey.com/en/products/detail/w%C3%BCrth-elektronik/7443551131/1638545"
    "G19" = "https://www.digikey.com/en/products/detail/bourns-inc/3352T-1-203LF/1088346"
    "G15" = "https://www.digikey.com/en/products/detail/texas-instruments/OPA4205APWR/17394950"
    "G20" = "https://www.digikey.com/en/products/detail/texas-instruments/SN75468DR/2255090?utm_adgroup=General&utm_source=google&utm_medium=cpc&utm_campaign=PMax%20Shopping_Product_Zombie%20SKUs&utm_term=&utm_content=General&utm_id=go_cmp-17815035045_adg-_ad-__dev-c_ext-_prd-2255090_sig-CjwKCAiAlJKuBhAdEiwAnZb7lY7edhjVnlVUhEyNawogcHzVo6bbfQ1LOtrzO4xh_eCL0cFOX98QUxoCbYMQAvD_BwE&gad_source=1&gclid=CjwKCAiAlJKuBhAdEiwAnZb7lY7edhjVnlVUhEyNawogcHzVo6bbfQ1LOtrzO4xh_eCL0cFOX98QUxoCbYMQAvD_BwE"
    "G21" = "https://www.digikey.com/en/products/detail/broadcom-limited/ACSL-6400-00TE/825239"
    "G22" = "https://www.digikey.com/en/products/detail/diodes-incorporated/1N4148WS-13-F/4249326?s=N4IgTCBcDaIIwDkAsckA4DqBlAtHAzDgGIgC6AvkA"
    "G14" = "https://www.digikey.com/en/products/detail/cui-devices/PJ-025/724801"
    "G23" = "https://www.digikey.com/en/products/detail/samtec-inc/TST-105-01-F-D/9497108"
    "G24" = "https://www.digikey.com/en/products/detail/samtec-inc/SSW-106-02-TM-S-RA/7891818"
    "G25" = "https://www.digikey.com/en/products/detail/e-switch/TL6330AF200Q/8032037"
    "G26" = "https://www.digikey.com/en/products/detail/murata-electronics/BLM21PG600SN1D/584263"
    "G27" = "https://www.digikey.com/en/products/detail/dialight/5988191107F/1291280"
    "G28" = "https://www.digikey.com/en/products/detail/liteon/LTW-170TK/758704"
    "G29" = "https://www.digikey.com/en/products/detail/dialight/5988110107F/1291272"
    "G30" = "https://www.digikey.com/en/products/detail/dialight/5988170107F/1291278"
}

# The text shown in column G mirrors the hyperlink target, so update it too.
foreach ($ref in $links.Keys) {
    $ws.Range($ref).Value = $links[$ref]
}

$ws.Hyperlinks.Delete()
foreach ($ref in $links.Keys) {
    $null = $ws.Hyperlinks.Add($ws.Range($ref), $links[$ref])
}

# ---------------------------------------------------------------------------
# Restore the selected cell as it was left in the saved workbook.
# ---------------------------------------------------------------------------
$null = $ws.Range("D7").Select()
